# -----------------------------------------------------------------------------
# Npnt-Itga8.xlsx - refresh with updated TPM-normalised expression values.
#
# Columns: E=Ligand-expressing cells, F=Ligand detection rate,
#          G/H=Ligand avg/total expr, I/J=Ligand derived specificity (avg/total),
#          M/N=Receptor avg/total expr, O/P=Receptor derived specificity (avg/total),
#          Q/R=Edge avg/total expression weight, S/T=Edge avg/total derived specificity.
# -----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Cells.Item(2, 5).Value = 3   # E2
$ws.Cells.Item(2, 6).Value = 1   # F2
$ws.Cells.Item(2, 7).Value = 1.660320666666667   # G2
$ws.Cells.Item(2, 8).Value = 4.980962   # H2
$ws.Cells.Item(2, 9).Value = 0.3342448133445559   # I2
$ws.Cells.Item(2, 10).Value = 0.3342448133445559   # J2
$ws.Cells.Item(2, 13).Value = 1.115279333333333   # M2
$ws.Cells.Item(2, 14).Value = 3.345838   # N2
$ws.Cells.Item(2, 15).Value = 0.08670767158519405   # O2
$ws.Cells.Item(2, 16).Value = 0.08670767158519403   # P2
$ws.Cells.Item(2, 17).Value = 1.851721326239555   # Q2
$ws.Cells.Item(2, 18).Value = 16.665491936156   # R2
$ws.Cells.Item(2, 19).Value = 0.02898158950453423   # S2
$ws.Cells.Item(2, 20).Value = 0.02898158950453423   # T2

# Row 3: ECs -> FAPs
$ws.Cells.Item(3, 5).Value = 3   # E3
$ws.Cells.Item(3, 6).Value = 1   # F3
$ws.Cells.Item(3, 7).Value = 1.660320666666667   # G3
$ws.Cells.Item(3, 8).Value = 4.980962   # H3
$ws.Cells.Item(3, 9).Value = 0.3342448133445559   # I3
$ws.Cells.Item(3, 10).Value = 0.3342448133445559   # J3
$ws.Cells.Item(3, 13).Value = 3.484068333333333   # M3
$ws.Cells.Item(3, 15).Value = 0.2708697667015328   # O3
$ws.Cells.Item(3, 16).Value = 0.2708697667015328   # P3
$ws.Cells.Item(3, 17).Value = 5.784670657912222   # Q3
$ws.Cells.Item(3, 18).Value = 52.06203592121   # R3
$ws.Cells.Item(3, 19).Value = 0.09053681461183721   # S3
$ws.Cells.Item(3, 20).Value = 0.09053681461183721   # T3

# Row 4: ECs -> MuSCs
$ws.Cells.Item(4, 5).Value = 3   # E4
$ws.Cells.Item(4, 6).Value = 1   # F4
$ws.Cells.Item(4, 7).Value = 1.660320666666667   # G4
$ws.Cells.Item(4, 8).Value = 4.980962   # H4
$ws.Cells.Item(4, 9).Value = 0.3342448133445559   # I4
$ws.Cells.Item(4, 10).Value = 0.3342448133445559   # J4
$ws.Cells.Item(4, 13).Value = 8.022733000000001   # M4
$ws.Cells.Item(4, 14).Value = 24.068199   # N4
$ws.Cells.Item(4, 15).Value = 0.6237293899283516   # O4
$ws.Cells.Item(4, 16).Value = 0.6237293899283515   # P4
$ws.Cells.Item(4, 17).Value = 13.32030940304867   # Q4
$ws.Cells.Item(4, 18).Value = 119.882784627438   # R4
$ws.Cells.Item(4, 19).Value = 0.2084783135141156   # S4
$ws.Cells.Item(4, 20).Value = 0.2084783135141155   # T4

# Row 5: ECs -> Resolving-Mac
$ws.Cells.Item(5, 5).Value = 3   # E5
$ws.Cells.Item(5, 6).Value = 1   # F5
$ws.Cells.Item(5, 7).Value = 1.660320666666667   # G5
$ws.Cells.Item(5, 8).Value = 4.980962   # H5
$ws.Cells.Item(5, 9).Value = 0.3342448133445559   # I5
$ws.Cells.Item(5, 10).Value = 0.3342448133445559   # J5
$ws.Cells.Item(5, 13).Value = 0.2404413333333334   # M5
$ws.Cells.Item(5, 14).Value = 0.7213240000000001   # N5
$ws.Cells.Item(5, 15).Value = 0.0186931717849216   # O5
$ws.Cells.Item(5, 16).Value = 0.0186931717849216   # P5
$ws.Cells.Item(5, 17).Value = 0.3992097148542222   # Q5
$ws.Cells.Item(5, 18).Value = 3.592887433688   # R5
$ws.Cells.Item(5, 19).Value = 0.006248095714068838   # S5
$ws.Cells.Item(5, 20).Value = 0.006248095714068837   # T5

# Row 6: FAPs -> ECs
$ws.Cells.Item(6, 9).Value = 0.02135559483851209   # I6
$ws.Cells.Item(6, 10).Value = 0.02135559483851209   # J6
$ws.Cells.Item(6, 13).Value = 1.115279333333333   # M6
$ws.Cells.Item(6, 14).Value = 3.345838   # N6
$ws.Cells.Item(6, 15).Value = 0.08670767158519405   # O6
$ws.Cells.Item(6, 16).Value = 0.08670767158519403   # P6
$ws.Cells.Item(6, 17).Value = 0.1183103187191111   # Q6
$ws.Cells.Item(6, 18).Value = 1.064792868472   # R6
$ws.Cells.Item(6, 19).Value = 0.001851693903764171   # S6
$ws.Cells.Item(6, 20).Value = 0.001851693903764171   # T6

# Row 7: FAPs -> FAPs
$ws.Cells.Item(7, 9).Value = 0.02135559483851209   # I7
$ws.Cells.Item(7, 10).Value = 0.02135559483851209   # J7
$ws.Cells.Item(7, 13).Value = 3.484068333333333   # M7
$ws.Cells.Item(7, 15).Value = 0.2708697667015328   # O7
$ws.Cells.Item(7, 16).Value = 0.2708697667015328   # P7
$ws.Cells.Item(7, 17).Value = 0.3695946142244445   # Q7
$ws.Cells.Item(7, 19).Value = 0.005784584991680227   # S7
$ws.Cells.Item(7, 20).Value = 0.005784584991680227   # T7

# Row 8: FAPs -> MuSCs
$ws.Cells.Item(8, 9).Value = 0.02135559483851209   # I8
$ws.Cells.Item(8, 10).Value = 0.02135559483851209   # J8
$ws.Cells.Item(8, 13).Value = 8.022733000000001   # M8
$ws.Cells.Item(8, 14).Value = 24.068199   # N8
$ws.Cells.Item(8, 15).Value = 0.6237293899283516   # O8
$ws.Cells.Item(8, 16).Value = 0.6237293899283515   # P8
$ws.Cells.Item(8, 17).Value = 0.8510622136173335   # Q8
$ws.Cells.Item(8, 18).Value = 7.659559922556   # R8
$ws.Cells.Item(8, 19).Value = 0.0133201121401822   # S8
$ws.Cells.Item(8, 20).Value = 0.0133201121401822   # T8

# Row 9: FAPs -> Resolving-Mac
$ws.Cells.Item(9, 9).Value = 0.02135559483851209   # I9
$ws.Cells.Item(9, 10).Value = 0.02135559483851209   # J9
$ws.Cells.Item(9, 13).Value = 0.2404413333333334   # M9
$ws.Cells.Item(9, 14).Value = 0.7213240000000001   # N9
$ws.Cells.Item(9, 15).Value = 0.0186931717849216   # O9
$ws.Cells.Item(9, 16).Value = 0.0186931717849216   # P9
$ws.Cells.Item(9, 17).Value = 0.02550633722844445   # Q9
$ws.Cells.Item(9, 18).Value = 0.229557035056   # R9
$ws.Cells.Item(9, 19).Value = 0.0003992038028854915   # S9
$ws.Cells.Item(9, 20).Value = 0.0003992038028854915   # T9

# Row 10: MuSCs -> ECs
$ws.Cells.Item(10, 7).Value = 3.108009   # G10
$ws.Cells.Item(10, 8).Value = 9.324027000000001   # H10
$ws.Cells.Item(10, 9).Value = 0.6256838868143543   # I10
$ws.Cells.Item(10, 10).Value = 0.6256838868143542   # J10
$ws.Cells.Item(10, 13).Value = 1.115279333333333   # M10
$ws.Cells.Item(10, 14).Value = 3.345838   # N10
$ws.Cells.Item(10, 15).Value = 0.08670767158519405   # O10
$ws.Cells.Item(10, 16).Value = 0.08670767158519403   # P10
$ws.Cells.Item(10, 17).Value = 3.466298205514001   # Q10
$ws.Cells.Item(10, 18).Value = 31.196683849626   # R10
$ws.Cells.Item(10, 19).Value = 0.05425159297404675   # S10
$ws.Cells.Item(10, 20).Value = 0.05425159297404673   # T10

# Row 11: MuSCs -> FAPs
$ws.Cells.Item(11, 7).Value = 3.108009   # G11
$ws.Cells.Item(11, 8).Value = 9.324027000000001   # H11
$ws.Cells.Item(11, 9).Value = 0.6256838868143543   # I11
$ws.Cells.Item(11, 10).Value = 0.6256838868143542   # J11
$ws.Cells.Item(11, 13).Value = 3.484068333333333   # M11
$ws.Cells.Item(11, 15).Value = 0.2708697667015328   # O11
$ws.Cells.Item(11, 16).Value = 0.2708697667015328   # P11
$ws.Cells.Item(11, 17).Value = 10.828515736615   # Q11
$ws.Cells.Item(11, 18).Value = 97.456641629535   # R11
$ws.Cells.Item(11, 19).Value = 0.1694788484503124   # S11
$ws.Cells.Item(11, 20).Value = 0.1694788484503124   # T11

# Row 12: MuSCs -> MuSCs
$ws.Cells.Item(12, 7).Value = 3.108009   # G12
$ws.Cells.Item(12, 8).Value = 9.324027000000001   # H12
$ws.Cells.Item(12, 9).Value = 0.6256838868143543   # I12
$ws.Cells.Item(12, 10).Value = 0.6256838868143542   # J12
$ws.Cells.Item(12, 13).Value = 8.022733000000001   # M12
$ws.Cells.Item(12, 14).Value = 24.068199   # N12
$ws.Cells.Item(12, 15).Value = 0.6237293899283516   # O12
$ws.Cells.Item(12, 16).Value = 0.6237293899283515   # P12
$ws.Cells.Item(12, 17).Value = 24.934726368597   # Q12
$ws.Cells.Item(12, 18).Value = 224.412537317373   # R12
$ws.Cells.Item(12, 19).Value = 0.390257429010717   # S12
$ws.Cells.Item(12, 20).Value = 0.3902574290107169   # T12

# Row 13: MuSCs -> Resolving-Mac
$ws.Cells.Item(13, 7).Value = 3.108009   # G13
$ws.Cells.Item(13, 8).Value = 9.324027000000001   # H13
$ws.Cells.Item(13, 9).Value = 0.6256838868143543   # I13
$ws.Cells.Item(13, 10).Value = 0.6256838868143542   # J13
$ws.Cells.Item(13, 13).Value = 0.2404413333333334   # M13
$ws.Cells.Item(13, 14).Value = 0.7213240000000001   # N13
$ws.Cells.Item(13, 15).Value = 0.0186931717849216   # O13
$ws.Cells.Item(13, 16).Value = 0.0186931717849216   # P13
$ws.Cells.Item(13, 17).Value = 0.7472938279720002   # Q13
$ws.Cells.Item(13, 18).Value = 6.725644451748002   # R13
$ws.Cells.Item(13, 19).Value = 0.01169601637927817   # S13
$ws.Cells.Item(13, 20).Value = 0.01169601637927816   # T13

# Row 14: Resolving-Mac -> ECs
$ws.Cells.Item(14, 7).Value = 0.09296800000000001   # G14
$ws.Cells.Item(14, 8).Value = 0.278904   # H14
$ws.Cells.Item(14, 9).Value = 0.01871570500257782   # I14
$ws.Cells.Item(14, 10).Value = 0.01871570500257782   # J14
$ws.Cells.Item(14, 13).Value = 1.115279333333333   # M14
$ws.Cells.Item(14, 14).Value = 3.345838   # N14
$ws.Cells.Item(14, 15).Value = 0.08670767158519405   # O14
$ws.Cells.Item(14, 16).Value = 0.08670767158519403   # P14
$ws.Cells.Item(14, 17).Value = 0.1036852890613333   # Q14
$ws.Cells.Item(14, 18).Value = 0.9331676015520001   # R14
$ws.Cells.Item(14, 19).Value = 0.001622795202848891   # S14
$ws.Cells.Item(14, 20).Value = 0.00162279520284889   # T14

# Row 15: Resolving-Mac -> FAPs
$ws.Cells.Item(15, 7).Value = 0.09296800000000001   # G15
$ws.Cells.Item(15, 8).Value = 0.278904   # H15
$ws.Cells.Item(15, 9).Value = 0.01871570500257782   # I15
$ws.Cells.Item(15, 10).Value = 0.01871570500257782   # J15
$ws.Cells.Item(15, 13).Value = 3.484068333333333   # M15
$ws.Cells.Item(15, 15).Value = 0.2708697667015328   # O15
$ws.Cells.Item(15, 16).Value = 0.2708697667015328   # P15
$ws.Cells.Item(15, 17).Value = 0.3239068648133334   # Q15
$ws.Cells.Item(15, 18).Value = 2.91516178332   # R15
$ws.Cells.Item(15, 19).Value = 0.005069518647702964   # S15
$ws.Cells.Item(15, 20).Value = 0.005069518647702964   # T15

# Row 16: Resolving-Mac -> MuSCs
$ws.Cells.Item(16, 7).Value = 0.09296800000000001   # G16
$ws.Cells.Item(16, 8).Value = 0.278904   # H16
$ws.Cells.Item(16, 9).Value = 0.01871570500257782   # I16
$ws.Cells.Item(16, 10).Value = 0.01871570500257782   # J16
$ws.Cells.Item(16, 13).Value = 8.022733000000001   # M16
$ws.Cells.Item(16, 14).Value = 24.068199   # N16
$ws.Cells.Item(16, 15).Value = 0.6237293899283516   # O16
$ws.Cells.Item(16, 16).Value = 0.6237293899283515   # P16
$ws.Cells.Item(16, 17).Value = 0.7458574415440001   # Q16
$ws.Cells.Item(16, 18).Value = 6.712716973896001   # R16
$ws.Cells.Item(16, 19).Value = 0.01167353526333686   # S16
$ws.Cells.Item(16, 20).Value = 0.01167353526333686   # T16

# Row 17: Resolving-Mac -> Resolving-Mac
$ws.Cells.Item(17, 7).Value = 0.09296800000000001   # G17
$ws.Cells.Item(17, 8).Value = 0.278904   # H17
$ws.Cells.Item(17, 9).Value = 0.01871570500257782   # I17
$ws.Cells.Item(17, 10).Value = 0.01871570500257782   # J17
$ws.Cells.Item(17, 13).Value = 0.2404413333333334   # M17
$ws.Cells.Item(17, 14).Value = 0.7213240000000001   # N17
$ws.Cells.Item(17, 15).Value = 0.0186931717849216   # O17
$ws.Cells.Item(17, 16).Value = 0.0186931717849216   # P17
$ws.Cells.Item(17, 17).Value = 0.02235334987733334   # Q17
$ws.Cells.Item(17, 18).Value = 0.2011801488960001   # R17
$ws.Cells.Item(17, 19).Value = 0.0003498558886891037   # S17
$ws.Cells.Item(17, 20).Value = 0.0003498558886891037   # T17
